$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9041280746459961
$ws.Range("B1").Value = 1.698358178138733
$ws.Range("C1").Value = 4.342281341552734
$ws.Range("D1").Value = 2.978245258331299
$ws.Range("E1").Value = 0.5283113718032837
